$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = '51.945.76'
$ws.Range("E2").Value = '  -0.68%  '
$ws.Range("D3").Value = '2.788.63'
$ws.Range("E3").Value = '  -2.15%  '
$ws.Range("E4").Value = '  +0.03%  '
Set-TextValue $ws "D5" '359.98'
$ws.Range("E5").Value = '  -0.22%  '
Set-TextValue $ws "D6" '109.60'
$ws.Range("E6").Value = '  -3.69%  '
Set-TextValue $ws "D7" '0.558'
$ws.Range("E7").Value = '  -2.92%  '
Set-TextValue $ws "D8" '0.999'
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -1.86%  '
Set-TextValue $ws "D10" '40.14'
$ws.Range("E10").Value = '  -3.56%  '
$ws.Range("E11").Value = '  -1.78%  '
$ws.Range("E12").Value = '  +1.12%  '
Set-TextValue $ws "D13" '19.52'
$ws.Range("E13").Value = '  -2.48%  '
$ws.Range("E14").Value = '  -2.86%  '
$ws.Range("D15").Value = '3.230.28'
$ws.Range("E15").Value = '  -1.90%  '
$ws.Range("D16").Value = '2.788.95'
$ws.Range("E16").Value = '  -1.47%  '
Set-TextValue $ws "D17" '0.937'
$ws.Range("E17").Value = '  +3.56%  '
$ws.Range("D18").Value = '51.897.40'
$ws.Range("E18").Value = '  -0.29%  '
$ws.Range("E19").Value = '  -1.42%  '
$ws.Range("E20").Value = '  -2.46%  '
Set-TextValue $ws "D21" '13.11'
$ws.Range("E21").Value = '  -3.73%  '
$ws.Range("E22").Value = '  -1.81%  '
Set-TextValue $ws "D23" '70.37'
$ws.Range("E23").Value = '  -0.06%  '
Set-TextValue $ws "D24" '269.79'
$ws.Range("E24").Value = '  +0.62%  '
Set-TextValue $ws "D25" '2.75'
$ws.Range("E25").Value = '  -2.54%  '
Set-TextValue $ws "D26" '26.54'
$ws.Range("E26").Value = '  -2.72%  '
$ws.Range("E27").Value = '  -0.03%  '
Set-TextValue $ws "D28" '0.161'
$ws.Range("E28").Value = '  +14.71%  '
Set-TextValue $ws "D29" '10.30'
$ws.Range("E29").Value = '  -1.38%  '
Set-TextValue $ws "D30" '2.28'
$ws.Range("E30").Value = '  +1.20%  '
$ws.Range("E31").Value = '  +2.06%  '
Set-TextValue $ws "D32" '51.97'
$ws.Range("E32").Value = '  -3.11%  '
Set-TextValue $ws "D33" '34.35'
$ws.Range("E33").Value = '  +0.15%  '
Set-TextValue $ws "D35" '0.0846'
$ws.Range("E35").Value = '  +0.06%  '
Set-TextValue $ws "D36" '5.24'
$ws.Range("E36").Value = '  -3.06%  '
$ws.Range("E37").Value = '  +0.16%  '
Set-TextValue $ws "D38" '19.04'
$ws.Range("E38").Value = '  +3.71%  '
$ws.Range("E39").Value = '  -2.27%  '
$ws.Range("E40").Value = '  -4.12%  '
Set-TextValue $ws "D41" '2.62'
$ws.Range("E41").Value = '  +1.81%  '
$ws.Range("E42").Value = '  -2.12%  '
$ws.Range("E43").Value = '  -1.26%  '
$ws.Range("B44").Value = 'Monero'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws "D44" '119.33'
$ws.Range("E44").Value = '  -7.09%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws "D45" '21.82'
$ws.Range("E45").Value = '  -8.41%  '
$ws.Range("D46").Value = '2.083.52'
$ws.Range("E46").Value = '  -1.62%  '
$ws.Range("E47").Value = '  -4.57%  '
$ws.Range("E49").Value = '  -1.36%  '
Set-TextValue $ws "D50" '0.958'
$ws.Range("E50").Value = '  -5.58%  '
$ws.Range("E51").Value = '  -2.42%  '
